# Add 2023 working-hours data (rows 32-43) and center the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fill in the 2023 monthly data (values + the "month abbreviation"
#    formula) first -- the formatting is applied afterwards by copying the
#    existing row's style so it matches the rest of the table exactly
#    (quote-prefixed text style for B/C, plain number style for D/E, etc.)
$dates = @(44927, 44958, 44986, 45017, 45047, 45078, 45108, 45139, 45170, 45200, 45231, 45261)
$labels = @("01. Янв 2023", "02. Фев 2023", "03. Мар 2023", "04. Апр 2023", "05. Май 2023", "06. Июн 2023", "07. Июл 2023", "08. Авг 2023", "09. Сен 2023", "10. Окт 2023", "11. Ноя 2023", "12. Дек 2023")
$hours = @(136, 143, 175, 160, 160, 168, 168, 184, 168, 176, 167, 168)
$hours90 = @(122.4, 128.6, 157.4, 144, 144, 151.19999999999999, 151.19999999999999, 165.6, 151.19999999999999, 158.4, 150.19999999999999, 151.19999999999999)

for ($i = 0; $i -lt 12; $i++) {
    $r = 32 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $labels[$i]
    $ws.Cells.Item($r, 4).Value = $hours[$i]
    $ws.Cells.Item($r, 5).Value = $hours90[$i]
}

# Fill the "month abbreviation" formula down column C as a single range
# assignment so it becomes one shared formula (C32:C43), matching the
# pattern already used for rows 3-31 (si="0").
$ws.Range("C32:C43").Formula = '=TEXT(A32,"МММ")'

# 2) Copy the cell formatting from the last populated data row (31) down
#    onto the twelve new rows so they pick up the existing date / text /
#    number styles instead of the old empty "placeholder" style. Doing this
#    last means the copied format "wins" over whatever default formatting
#    the value/formula assignments above may have picked.
$ws.Range("A31:E31").Copy()
$ws.Range("A32:E43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Center the header row text.
$ws.Range("A1:E1").HorizontalAlignment = -4108

# 4) The month-abbreviation column is now only 3 Cyrillic letters long (vs.
#    the old placeholder text), so re-fit columns C and D to the new
#    narrower content.
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
